$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.885185185187
$ws.Range("B3").Value = -28
$ws.Range("C3").Value = 52
$ws.Range("D3").Value = 47
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 81
$ws.Range("G3").Value = 12244
$ws.Range("H3").Value = 10104
$ws.Range("I3").Value = 1175
$ws.Range("J3").Value = 131
$ws.Range("K3").Value = 119
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 26
$ws.Range("N3").Value = "Bag"
